$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated P&L calculation: the flight revenue summary table now covers
# 25 flights (rows 2-26) instead of the previous 10 (rows 2-11), with
# refreshed revenue figures for every flight.

$row1 = ,@("FL2240", 36651.01, 23424.5, 4731.21, 6750, 1745.3, 89, 411.81)
$row2 = ,@("FL8121", 38517.45, 19014.94, 7118.35, 10550, 1834.16, 138, 279.11)
$row3 = ,@("FL9519", 196959.99, 166232.53, 8548.41, 12800, 9379.05, 168, 1172.38)
$row4 = ,@("FL3335", 42465.4, 23234.73, 6808.51, 10400, 2022.16, 131, 324.16)
$row5 = ,@("FL1198", 36529.58, 19102.98, 5637.09, 10050, 1739.51, 122, 299.42)
$row6 = ,@("FL1294", 39121.91, 25581.72, 4127.24, 7550, 1862.95, 96, 407.52)
$row7 = ,@("FL8464", 38593.74, 23556.54, 5049.39, 8150, 1837.81, 106, 364.09)
$row8 = ,@("FL1344", 44625.37, 24852.39, 6697.97, 10950, 2125.01, 142, 314.26)
$row9 = ,@("FL2664", 39940.18, 21827.41, 6110.84, 10100, 1901.93, 124, 322.1)
$row10 = ,@("FL3756", 39563.61, 19541.63, 7738, 10400, 1883.98, 146, 270.98)
$row11 = ,@("FL5997", 41483.9, 27764.05, 4944.43, 6800, 1975.42, 99, 419.03)
$row12 = ,@("FL9566", 38059.08, 25143.4, 4653.34, 6450, 1812.34, 94, 404.88)
$row13 = ,@("FL1990", 114938.91, 87579.06, 9486.56, 12400, 5473.29, 183, 628.08)
$row14 = ,@("FL8118", 39873.34, 24109.15, 5165.46, 8700, 1898.73, 107, 372.65)
$row15 = ,@("FL2977", 230649.16, 196444.7, 9821.16, 13400, 10983.3, 200, 1153.25)
$row16 = ,@("FL5421", 42894.38, 27943.29, 4858.51, 8050, 2042.58, 103, 416.45)
$row17 = ,@("FL5301", 50783.06, 27951.26, 8363.56, 12050, 2418.24, 160, 317.39)
$row18 = ,@("FL9612", 131931.63, 99149.57, 10249.58, 16250, 6282.48, 203, 649.91)
$row19 = ,@("FL6575", 48544.14, 26800.73, 7481.79, 11950, 2311.62, 154, 315.22)
$row20 = ,@("FL4876", 32434.01, 21090.71, 3998.83, 5800, 1544.47, 78, 415.82)
$row21 = ,@("FL1573", 39805.96, 26148.53, 4611.91, 7150, 1895.52, 96, 414.65)
$row22 = ,@("FL6970", 211855.36, 179618.13, 9298.88, 12850, 10088.35, 183, 1157.68)
$row23 = ,@("FL5705", 45463.11, 23939.91, 7858.3, 11500, 2164.9, 151, 301.08)
$row24 = ,@("FL1013", 101040.83, 75832.24, 8047.11, 12350, 4811.48, 161, 627.58)
$row25 = ,@("FL8222", 39099.7, 21096.29, 6791.52, 9350, 1861.89, 133, 293.98)

$data = $row1 + $row2 + $row3 + $row4 + $row5 + $row6 + $row7 + $row8 + $row9 + $row10 + $row11 + $row12 + $row13 + $row14 + $row15 + $row16 + $row17 + $row18 + $row19 + $row20 + $row21 + $row22 + $row23 + $row24 + $row25

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    $r = $startRow + $i
    $arr = New-Object 'object[,]' 1,8
    for ($c = 0; $c -lt 8; $c++) {
        $arr[0,$c] = $rowVals[$c]
    }
    $targetRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 8))
    $targetRange.Value = $arr
}
